# edit.ps1 - apply the edits described in the commit diff to the active document
$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2  (used throughout below)

# ---------------------------------------------------------------------------
# 1) "... characters were embedded in the elevation data." ->
#    "... characters were embedded in some fields."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "embedded in the elevation data.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "embedded in some fields.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "... retained only one potential issue of a city having multiple
#    associated counties, which **could** impact its ability to merge
#    efficiently with other tables." ->
#    "... retained only a few missing values which were manually entered
#    into the exported Excel file, where final formatting of the data
#    fields was also done."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "one potential issue of a city having multiple associated counties, which **could** impact its ability to merge efficiently with other tables.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "a few missing values which were manually entered into the exported Excel file, where final formatting of the data fields was also done.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Tech used: Webscaping, Python, ..." -> "Tech used: Webscaping, Excel, Python, ..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Webscaping, Python",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Webscaping, Excel, Python",
    2) | Out-Null

Write-Output "done"
